# "Implemented getting kafka relations."
#
# The structure-extraction run that produced this workbook re-walked the
# Java reflection metadata after the kafka-relations feature was added, and
# the (non-deterministic) field enumeration order shifted for several
# already-documented classes/enums on the "classFields" sheet even though
# the set of fields per class is unchanged. Net effect: a bunch of
# Field Name / Field Modifier / Field Type cells on existing rows 9-98 now
# hold a different value than before. Apply each corrected cell in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# org.andante.activity.enums.Priority (rows 9-13)
$ws.Cells.Item(9, 3).Value = 'private'
$ws.Cells.Item(9, 4).Value = 'org.andante.activity.enums.Priority[]'
$ws.Cells.Item(13, 3).Value = 'public'
$ws.Cells.Item(13, 4).Value = 'org.andante.activity.enums.Priority'

# org.andante.activity.dto.NewsletterDTO$NewsletterDTOBuilder (rows 14-16)
$ws.Cells.Item(15, 4).Value = 'java.time.LocalDateTime'
$ws.Cells.Item(16, 4).Value = 'java.lang.String'

# org.andante.activity.dto.UserProfileDTO$UserProfileDTOBuilder (rows 20-24)
$ws.Cells.Item(20, 4).Value = 'java.util.Set'
$ws.Cells.Item(22, 4).Value = 'java.lang.String'

# org.andante.activity.dto.ActivityQuerySpecification (rows 25-27)
$ws.Cells.Item(25, 4).Value = 'java.lang.String'
$ws.Cells.Item(27, 4).Value = 'java.lang.Integer'

# org.andante.activity.dto.UserImageDTO$UserImageDTOBuilder (rows 31-32)
$ws.Cells.Item(31, 2).Value = 'imageUrl'
$ws.Cells.Item(32, 2).Value = 'username'

# org.andante.activity.dto.ActivityQuerySpecification$ActivityQuerySpecificationBuilder (rows 44-46)
$ws.Cells.Item(44, 2).Value = 'query'
$ws.Cells.Item(45, 2).Value = 'pageSize'
$ws.Cells.Item(45, 4).Value = 'java.lang.String'
$ws.Cells.Item(46, 2).Value = 'pageNumber'
$ws.Cells.Item(46, 4).Value = 'java.lang.Integer'

# org.andante.activity.enums.Domain (rows 51-53)
$ws.Cells.Item(51, 2).Value = 'HIGHEST'
$ws.Cells.Item(51, 3).Value = 'private'
$ws.Cells.Item(51, 4).Value = 'org.andante.activity.enums.Domain[]'
$ws.Cells.Item(52, 2).Value = 'ORDER'
$ws.Cells.Item(53, 2).Value = 'PRODUCT'
$ws.Cells.Item(53, 3).Value = 'public'
$ws.Cells.Item(53, 4).Value = 'org.andante.activity.enums.Domain'

# org.andante.activity.dto.ActivityDTO$ActivityDTOBuilder (rows 57-64)
$ws.Cells.Item(57, 4).Value = 'org.andante.activity.enums.Priority'
$ws.Cells.Item(58, 4).Value = 'org.andante.activity.enums.Domain'
$ws.Cells.Item(59, 4).Value = 'java.util.Set'
$ws.Cells.Item(60, 2).Value = 'street'
$ws.Cells.Item(61, 2).Value = 'relatedId'
$ws.Cells.Item(61, 4).Value = 'java.lang.String'
$ws.Cells.Item(62, 2).Value = 'acknowledgedUsers'
$ws.Cells.Item(63, 4).Value = 'java.time.LocalDateTime'
$ws.Cells.Item(64, 4).Value = 'java.util.Set'

# org.andante.activity.dto.ActivityDTO (rows 68-75)
$ws.Cells.Item(68, 2).Value = 'affectedUsers'
$ws.Cells.Item(68, 4).Value = 'org.andante.activity.enums.Domain'
$ws.Cells.Item(69, 2).Value = 'relatedId'
$ws.Cells.Item(70, 4).Value = 'java.time.LocalDateTime'
$ws.Cells.Item(71, 2).Value = 'id'
$ws.Cells.Item(71, 4).Value = 'org.andante.activity.enums.Priority'
$ws.Cells.Item(72, 2).Value = 'eventTimestamp'
$ws.Cells.Item(73, 2).Value = 'domain'
$ws.Cells.Item(73, 4).Value = 'java.util.Set'
$ws.Cells.Item(74, 2).Value = 'street'
$ws.Cells.Item(75, 4).Value = 'java.lang.String'

# org.andante.activity.dto.UserImageDTO$UserImageDTOBuilder (rows 76-77)
$ws.Cells.Item(76, 2).Value = 'imageUrl'
$ws.Cells.Item(77, 2).Value = 'username'

# org.andante.activity.dto.NewsletterDTO (rows 78-80)
$ws.Cells.Item(78, 2).Value = 'isConfirmed'
$ws.Cells.Item(78, 4).Value = 'java.lang.Boolean'
$ws.Cells.Item(79, 4).Value = 'java.lang.String'
$ws.Cells.Item(80, 2).Value = 'emailAddress'
$ws.Cells.Item(80, 4).Value = 'java.time.LocalDateTime'

# org.andante.activity.dto.UpdatedUserDetailsDTO (rows 81-91)
$ws.Cells.Item(81, 2).Value = 'dateOfBirth'
$ws.Cells.Item(82, 2).Value = 'description'
$ws.Cells.Item(83, 2).Value = 'gender'
$ws.Cells.Item(87, 2).Value = 'lastName'
$ws.Cells.Item(88, 2).Value = 'city'
$ws.Cells.Item(89, 2).Value = 'country'
$ws.Cells.Item(90, 2).Value = 'postalCode'
$ws.Cells.Item(91, 2).Value = 'street'

# org.andante.activity.dto.UserProfileDTO (rows 92-96)
$ws.Cells.Item(92, 2).Value = 'username'
$ws.Cells.Item(92, 4).Value = 'java.lang.String'
$ws.Cells.Item(93, 2).Value = 'key'
$ws.Cells.Item(93, 4).Value = 'java.util.Set'
$ws.Cells.Item(94, 2).Value = 'observedUsers'
$ws.Cells.Item(95, 2).Value = 'observingUsers'
$ws.Cells.Item(95, 4).Value = 'java.util.Set'
$ws.Cells.Item(96, 2).Value = 'imageUrl'
$ws.Cells.Item(96, 4).Value = 'java.lang.String'

# org.andante.activity.exception.UserConflictException (rows 97-98)
$ws.Cells.Item(97, 2).Value = 'serialVersionUID'
$ws.Cells.Item(97, 4).Value = 'long'
$ws.Cells.Item(98, 2).Value = '$assertionsDisabled'
$ws.Cells.Item(98, 4).Value = 'boolean'
